# Applies the "Fixed a little error" edit:
#  1. Rebuild the paragraph-properties of the "Specificity" heading
#     paragraph (widowControl / tabs / autoSpaceDE / autoSpaceDN /
#     adjustRightInd / spacing-after) and relocate the "_GoBack"
#     bookmark from mid-sentence in the following paragraph to the
#     very start of the heading paragraph.
#  2. Merge several runs that were needlessly split across multiple
#     <w:r> elements (no visible text change, just tidies the XML) by
#     doing a no-op Find/Replace across each run boundary, which makes
#     Word re-flow the matched span into a single run.

$d = $word.ActiveDocument

# --- 1. "Specificity" heading paragraph: new pPr + moved bookmark ---------
$heading = $d.Paragraphs.Item(3)
$headingFrag = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:widowControl w:val="0"/>
    <w:tabs>
      <w:tab w:val="left" w:pos="220"/>
      <w:tab w:val="left" w:pos="720"/>
    </w:tabs>
    <w:autoSpaceDE w:val="0"/>
    <w:autoSpaceDN w:val="0"/>
    <w:adjustRightInd w:val="0"/>
    <w:spacing w:after="260" w:line="480" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Helvetica Neue"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Helvetica Neue"/>
    </w:rPr>
    <w:t>Specificity</w:t>
  </w:r>
</w:p>
'@
$heading.Range.InsertXML($headingFrag)

# Re-seat the lone "_GoBack" bookmark at the start of that paragraph;
# Bookmarks.Add with an existing name moves it (Word only allows one
# bookmark per name) instead of creating a duplicate, so the old
# occurrence (mid-way through the next paragraph's text) disappears.
$heading = $d.Paragraphs.Item(3)
$headingStart = $d.Range($heading.Range.Start, $heading.Range.Start)
$d.Bookmarks.Add("_GoBack", $headingStart)

# --- 2. Merge runs that were split without any formatting difference -----
# Each Find/Replace only needs to span the run boundary; Word collapses
# the matched text into a single run using the formatting already present.

$d.Content.Find.Execute("pints. This", $true, $false, $false, $false, $false, `
    $true, 1, $false, "pints. This", 2) | Out-Null

$d.Content.Find.Execute("Site point. 3 June. 2012.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Site point. 3 June. 2012.", 2) | Out-Null

$d.Content.Find.Execute("CSS-Tricks. 3 June. 2012", $true, $false, $false, $false, $false, `
    $true, 1, $false, "CSS-Tricks. 3 June. 2012", 2) | Out-Null

$d.Content.Find.Execute("http://css-tricks.com/specifics-on-css-specificity/", $true, $false, $false, $false, $false, `
    $true, 1, $false, "http://css-tricks.com/specifics-on-css-specificity/", 2) | Out-Null

Write-Output "Edit complete"
